$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.542.09'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.754.34'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4597'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.95%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3562'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07464'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.48'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.086'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.80%  '
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.76'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.54%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.011'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.179'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.752.74'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '94.27'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001056'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06408'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.0000'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.746'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.599.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.083'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.34'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.14'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.950.46'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.114'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.58'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.082'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09226'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.655'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.510'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.55%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02285'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.73%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '11.72'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2088'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06013'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6285'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.924'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.181'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.390'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.781'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.11'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5857'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.87'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.936'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06893'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.128'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.94%  '
